$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Browser column (F2) value from "Chrome" to "Firefox"
$ws.Range("F2").Value = "Firefox"

# Move the active selection to G2, matching the saved workbook state
$ws.Range("G2").Select()
